$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells with corrected values
$ws.Range("F545").Value = 16775

$ws.Range("F575").Value = 26636

$ws.Range("F606").Value = 14615

$ws.Range("F636").Value = 50196
$ws.Range("G636").Value = 2360

$ws.Range("F637").Value = 43653

$ws.Range("F638").Value = 37761

$ws.Range("F639").Value = 40731

$ws.Range("F642").Value = 67407

$ws.Range("F644").Value = 36898
$ws.Range("G644").Value = 1503

$ws.Range("F645").Value = 35497

$ws.Range("F646").Value = 36017

$ws.Range("F649").Value = 62390

$ws.Range("F650").Value = 38060
$ws.Range("G650").Value = 1186

$ws.Range("F651").Value = 37123
$ws.Range("G651").Value = 1061

$ws.Range("F652").Value = 34992
$ws.Range("G652").Value = 1093

$ws.Range("F653").Value = 34142
$ws.Range("G653").Value = 1012

$ws.Range("F654").Value = 14322
$ws.Range("G654").Value = 684

$ws.Range("F655").Value = 25279
$ws.Range("G655").Value = 805

$ws.Range("F656").Value = 52251
$ws.Range("G656").Value = 1236

$ws.Range("F657").Value = 34053
$ws.Range("G657").Value = 872

$ws.Range("F658").Value = 27150
$ws.Range("G658").Value = 785

$ws.Range("F659").Value = 26191
$ws.Range("G659").Value = 845

$ws.Range("F660").Value = 6175
$ws.Range("G660").Value = 252

$ws.Range("F661").Value = 4857

$ws.Range("F662").Value = 12593
$ws.Range("G662").Value = 544

$ws.Range("F663").Value = 37034
$ws.Range("G663").Value = 1154

$ws.Range("F664").Value = 26671
$ws.Range("G664").Value = 778

$ws.Range("F665").Value = 24289
$ws.Range("G665").Value = 656

$ws.Range("F666").Value = 23669
$ws.Range("G666").Value = 758

$ws.Range("F667").Value = 17045
$ws.Range("G667").Value = 598

$ws.Range("F668").Value = 3326
$ws.Range("G668").Value = 125

$ws.Range("F669").Value = 22913
$ws.Range("G669").Value = 622

$ws.Range("F670").Value = 51733
$ws.Range("G670").Value = 882

# Add new rows 671-676
$newRows = @(
    @(671, 44565, 853416, 13498, 3900, 16824, 32343, 605),
    @(672, 44566, 856685, 12432, 3269, 16871, 29201, 568),
    @(673, 44567, 857801, 5070, 1116, 16896, 9914, 304),
    @(674, 44568, 860832, 13258, 3031, 16933, 24539, 609),
    @(675, 44569, 862962, 9378, 2130, 16962, 11868, 314),
    @(676, 44570, 863817, 4490, 855, 16989, 19794, 297)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
